$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (recomputed TPM-based NATMI ligand-receptor values). The
# "Resolving-Mac" cluster has been dropped as a *target* cluster, so the
# table shrinks from 12 data rows (4 targets x 3 senders) down to 9 data
# rows (3 targets x 3 senders); it is still present as a *sending* cluster.
$data = @(
    ,@("ECs","Pspn","Ret","ECs",2,0.6666666666666666,1.863243333333333,5.589729999999999,0.6067417803684044,0.6067417803684044,3,1,4.132875333333334,12.398626,0.391502049440379,0.3915020494403789,7.700552412331111,69.30497171098,0.2375406504953346,0.2375406504953346)
    ,@("ECs","Pspn","Ret","FAPs",2,0.6666666666666666,1.863243333333333,5.589729999999999,0.6067417803684044,0.6067417803684044,3,1,5.094146333333334,15.282439,0.4825620346115429,0.4825620346115428,9.491634194607778,85.42470775147,0.2927905480184071,0.2927905480184071)
    ,@("ECs","Pspn","Ret","MuSCs",2,0.6666666666666666,1.863243333333333,5.589729999999999,0.6067417803684044,0.6067417803684044,3,1,1.329437333333333,3.988312,0.1259359159480782,0.1259359159480782,2.477065248417778,22.29358723576,0.07641058185466269,0.07641058185466269)
    ,@("MuSCs","Pspn","Ret","ECs",3,1,0.9841983333333334,2.952595,0.3204918210730839,0.3204918210730839,3,1,4.132875333333334,12.398626,0.391502049440379,0.3915020494403789,4.067569014941112,36.60812113447,0.1254732047789916,0.1254732047789916)
    ,@("MuSCs","Pspn","Ret","FAPs",3,1,0.9841983333333334,2.952595,0.3204918210730839,0.3204918210730839,3,1,5.094146333333334,15.282439,0.4825620346115429,0.4825620346115428,5.013650331022778,45.122852979205,0.154657185253386,0.1546571852533859)
    ,@("MuSCs","Pspn","Ret","MuSCs",3,1,0.9841983333333334,2.952595,0.3204918210730839,0.3204918210730839,3,1,1.329437333333333,3.988312,0.1259359159480782,0.1259359159480782,1.308430007737778,11.77587006964,0.04036143104070641,0.04036143104070641)
    ,@("Resolving-Mac","Pspn","Ret","ECs",1,0.3333333333333333,0.2234583333333333,0.6703750000000001,0.07276639855851162,0.07276639855851162,3,1,4.132875333333334,12.398626,0.391502049440379,0.3915020494403789,0.9235254338611112,8.311728904750002,0.02848819416605274,0.02848819416605274)
    ,@("Resolving-Mac","Pspn","Ret","FAPs",1,0.3333333333333333,0.2234583333333333,0.6703750000000001,0.07276639855851162,0.07276639855851162,3,1,5.094146333333334,15.282439,0.4825620346115429,0.4825620346115428,1.138329449402778,10.244965044625,0.03511430133974981,0.03511430133974981)
    ,@("Resolving-Mac","Pspn","Ret","MuSCs",1,0.3333333333333333,0.2234583333333333,0.6703750000000001,0.07276639855851162,0.07276639855851162,3,1,1.329437333333333,3.988312,0.1259359159480782,0.1259359159480782,0.2970738507777778,2.673664657,0.009163903052709078,0.009163903052709078)
)

# Remove the three rows that no longer exist in the refreshed data (old
# rows 11-13) before overwriting the remaining rows, so nothing stale is
# left below row 10.
$ws.Range("A11:T13").EntireRow.Delete()

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowValues[$c]
    }
}
